# Updates cryptos list values (Price and Volume(1h) columns) to reflect
# the latest scrape, including the TheSandbox/Aptos row swap at rows 40-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.476.92'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.798.23'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5444'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3786'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07512'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.110'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.166'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.301'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').Value = '1.797.20'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001066'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06497'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.963'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '28.481.98'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.081'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').Value = '2.001.32'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.316'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.108'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1052'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.619'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.654'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2284'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06468'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02299'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.604'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.031'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.21'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.72%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6209'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.452'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.193'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.44'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5833'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.201'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.952'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06877'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.28%  '
